# The deck currently ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colour scheme (only wired to the notes master)
#   ppt/theme/theme2.xml -> "Integral" colour scheme (wired to the slide master / presentation)
#
# The authored edit swaps the two themes' content, so the presentation (and its
# slide master) now carries the "Office" 12-slot colour palette that used to
# live in theme1.xml, while theme1.xml (not reachable from the object model -
# it backs only the notes master) keeps its own palette unchanged.
#
# Font scheme and format scheme are identical between the two theme parts
# already, so the only observable difference is the 12-slot colour scheme.
# We rewrite those 12 slots on the presentation's live theme (reached through
# any slide's ThemeColorScheme - it is one shared part) to the target
# ("Office") values.

function HexToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette (was theme1.xml "Office" scheme), in clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$targetHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $targetHex.Count; $i++) {
    $tcs.Colors($i).RGB = HexToRgbLong $targetHex[$i - 1]
}
